# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table (rows 3-9) gets refreshed data and is re-sorted
# ascending by "Good Roaming Calculation (%)" (column D) ---

# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.0.3
$ws.Range("C3").Value = 3079
$ws.Range("D3").Value = 77.40000000000001

# Row 4: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.1.1
$ws.Range("C4").Value = 286
$ws.Range("D4").Value = 91.8

# Row 5: now Intel(R) Wi-Fi 6E AX211 160MHz - 22.170.0.3
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.170.0.3"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1950
$ws.Range("D5").Value = 96.5

# Row 6: now Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3
$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 135
$ws.Range("D6").Value = 96.7

# Row 7: now Intel(R) Wireless-AC 9560 160MHz - 20.110.0.3
$ws.Range("A7").Value = "Intel(R) Wireless-AC 9560 160MHz - 20.110.0.3"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 96.8

# Row 8: now Intel(R) Wi-Fi 6 AX201 160MHz - 22.120.0.3
$ws.Range("A8").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.120.0.3"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 839
$ws.Range("D8").Value = 98.09999999999999

# Row 9: Intel(R) Dual Band Wireless-AC 7260 - 18.33.17.1 (text/B/D unchanged)
$ws.Range("C9").Value = 154

# Row 10: Totals
$ws.Range("C10").Value = 6449

# --- "Good Drivers" table updates (Total Samples, column B) ---
$ws.Range("B20").Value = 56069
$ws.Range("B21").Value = 449371
$ws.Range("B25").Value = 276086
$ws.Range("B26").Value = 625298
$ws.Range("B31").Value = 453652
$ws.Range("B39").Value = 96091
$ws.Range("B42").Value = 99549
$ws.Range("B43").Value = 77999
$ws.Range("B47").Value = 175767
$ws.Range("B48").Value = 240182
$ws.Range("B54").Value = 684728
$ws.Range("B56").Value = 210188
$ws.Range("B60").Value = 308481
$ws.Range("B65").Value = 443223
$ws.Range("B66").Value = 109665
$ws.Range("B68").Value = 62515
